$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.4423115263591342
$ws.Range("C2").Value = 0.3792800776886729
$ws.Range("D2").Value = 0.04153582648314824
$ws.Range("F2").Value = 0.7749282079033861
$ws.Range("G2").Value = 0.6187996857315312
$ws.Range("H2").Value = 0.7493751252327172
$ws.Range("K2").Value = 0.2514287238144277
$ws.Range("L2").Value = 0.3031147674088714
$ws.Range("M2").Value = 0.1656964608235647
$ws.Range("O2").Value = 2.721674786174674
# Row 3
$ws.Range("B3").Value = 0.4042229608625973
$ws.Range("C3").Value = 0.381780301871423
$ws.Range("D3").Value = 0.03871279724742749
$ws.Range("F3").Value = 0.7771467947007764
$ws.Range("G3").Value = 0.6226241675176354
$ws.Range("H3").Value = 0.754765332532692
$ws.Range("K3").Value = 0.2196857830806493
$ws.Range("L3").Value = 0.2999378463029885
$ws.Range("M3").Value = 0.1579297760462772
$ws.Range("O3").Value = 2.740821529504998
# Row 4
$ws.Range("B4").Value = 0.3809050649576591
$ws.Range("C4").Value = 0.3834056892649116
$ws.Range("D4").Value = 0.03696360060506976
$ws.Range("F4").Value = 0.7789522475964858
$ws.Range("G4").Value = 0.6253550821754672
$ws.Range("H4").Value = 0.7583738548702428
$ws.Range("K4").Value = 0.2001117936538463
$ws.Range("L4").Value = 0.2981422160951368
$ws.Range("M4").Value = 0.1532242264532933
$ws.Range("O4").Value = 2.754006130871105
# Row 5
$ws.Range("B5").Value = 0.3714207485209897
$ws.Range("C5").Value = 0.3840907824670765
$ws.Range("D5").Value = 0.03624683910025794
$ws.Range("F5").Value = 0.7797995033374931
$ws.Range("G5").Value = 0.6265641624967344
$ws.Range("H5").Value = 0.7599195923842359
$ws.Range("K5").Value = 0.1921146889346801
$ws.Range("L5").Value = 0.2974495415190077
$ws.Range("M5").Value = 0.151322712420253
$ws.Range("O5").Value = 2.759738264361559
# Row 6
$ws.Range("B6").Value = 0.3698469863861931
$ws.Range("C6").Value = 0.3842059160614859
$ws.Range("D6").Value = 0.03612758387068027
$ws.Range("F6").Value = 0.7799469268820118
$ws.Range("G6").Value = 0.62677073980975
$ws.Range("H6").Value = 0.7601808070458915
$ws.Range("K6").Value = 0.1907855496858417
$ws.Range("L6").Value = 0.2973368848799112
$ws.Range("M6").Value = 0.1510079400072293
$ws.Range("O6").Value = 2.76071178498978
# Row 7
$ws.Range("B7").Value = 0.3807770827390868
$ws.Range("C7").Value = 0.3834148365621743
$ws.Range("D7").Value = 0.03695395005474467
$ws.Range("F7").Value = 0.7789632223616607
$ws.Range("G7").Value = 0.6253709987223388
$ws.Range("H7").Value = 0.7583943965079314
$ws.Range("K7").Value = 0.2000040244244303
$ws.Range("L7").Value = 0.2981327161867355
$ws.Range("M7").Value = 0.1531985168606447
$ws.Range("O7").Value = 2.754081981481889
# Row 8
$ws.Range("B8").Value = 0.4291647244641297
$ws.Range("C8").Value = 0.3801234517826195
$ws.Range("D8").Value = 0.04056575476774071
$ws.Range("F8").Value = 0.7756012214079533
$ws.Range("G8").Value = 0.6200389351644944
$ws.Range("H8").Value = 0.7511716719657642
$ws.Range("K8").Value = 0.2405014194834649
$ws.Range("L8").Value = 0.3019872386627753
$ws.Range("M8").Value = 0.1630054646976902
$ws.Range("O8").Value = 2.727980180185497
# Row 9
$ws.Range("B9").Value = 0.5245733200372058
$ws.Range("C9").Value = 0.374383090613069
$ws.Range("D9").Value = 0.0475215424716211
$ws.Range("F9").Value = 0.7725230952948721
$ws.Range("G9").Value = 0.6126201386812369
$ws.Range("H9").Value = 0.7393765984075742
$ws.Range("K9").Value = 0.3192342957632945
$ws.Range("L9").Value = 0.3107731131640605
$ws.Range("M9").Value = 0.1827336038928919
$ws.Range("O9").Value = 2.688124623594248
# Row 10
$ws.Range("B10").Value = 0.5949621389601987
$ws.Range("C10").Value = 0.3705981642000893
$ws.Range("D10").Value = 0.05255336505203445
$ws.Range("F10").Value = 0.7724025754784734
$ws.Range("G10").Value = 0.6090232965665905
$ws.Range("H10").Value = 0.7321507047506657
$ws.Range("K10").Value = 0.3766441734435944
$ws.Range("L10").Value = 0.3179735603485199
$ws.Range("M10").Value = 0.1975256737194329
$ws.Range("O10").Value = 2.66574612493946
# Row 11
$ws.Range("B11").Value = 0.6270423514181687
$ws.Range("C11").Value = 0.3689696327839336
$ws.Range("D11").Value = 0.05482517544919574
$ws.Range("F11").Value = 0.7728123460345913
$ws.Range("G11").Value = 0.607790008788939
$ws.Range("H11").Value = 0.7291752855341542
$ws.Range("K11").Value = 0.4026632008947502
$ws.Range("L11").Value = 0.3214106041756537
$ws.Range("M11").Value = 0.2043186660632657
$ws.Range("O11").Value = 2.657063910368521
# Row 12
$ws.Range("B12").Value = 0.6391983475205052
$ws.Range("C12").Value = 0.3683663165357132
$ws.Range("D12").Value = 0.05568294885313207
$ws.Range("F12").Value = 0.7730342777326555
$ws.Range("G12").Value = 0.6073809656939773
$ws.Range("H12").Value = 0.7280933241036252
$ws.Range("K12").Value = 0.412501544265865
$ws.Range("L12").Value = 0.3227352847823539
$ws.Range("M12").Value = 0.2069000819821127
$ws.Range("O12").Value = 2.653991504360732
# Row 13
$ws.Range("B13").Value = 0.6365799954937472
$ws.Range("C13").Value = 0.3684956572763625
$ws.Range("D13").Value = 0.05549832435465873
$ws.Range("F13").Value = 0.7729835123541022
$ws.Range("G13").Value = 0.6074664813909294
$ws.Range("H13").Value = 0.7283243537584525
$ws.Range("K13").Value = 0.4103833344714758
$ws.Range("L13").Value = 0.3224489629396743
$ws.Range("M13").Value = 0.2063437273954065
$ws.Range("O13").Value = 2.654643624301514
# Row 14
$ws.Range("B14").Value = 0.6280422779504704
$ws.Range("C14").Value = 0.3689197297779785
$ws.Range("D14").Value = 0.05489579555376167
$ws.Range("F14").Value = 0.7728292667874399
$ws.Range("G14").Value = 0.6077551944415802
$ws.Range("H14").Value = 0.7290853750588013
$ws.Range("K14").Value = 0.4034729010075182
$ws.Range("L14").Value = 0.3215191230135162
$ws.Range("M14").Value = 0.2045308600268072
$ws.Range("O14").Value = 2.656806825730484
# Row 15
$ws.Range("B15").Value = 0.6228136892607097
$ws.Range("C15").Value = 0.3691812269604853
$ws.Range("D15").Value = 0.05452640095927563
$ws.Range("F15").Value = 0.7727434796708863
$ws.Range("G15").Value = 0.6079395907544978
$ws.Range("H15").Value = 0.7295573507550017
$ws.Range("K15").Value = 0.3992381567438201
$ws.Range("L15").Value = 0.3209525812045086
$ws.Range("M15").Value = 0.2034216015684436
$ws.Range("O15").Value = 2.658159894327099
# Row 16
$ws.Range("B16").Value = 0.5928667679780801
$ws.Range("C16").Value = 0.3707064659039467
$ws.Range("D16").Value = 0.05240454756813051
$ws.Range("F16").Value = 0.7723851408008429
$ws.Range("G16").Value = 0.609112005879993
$ws.Range("H16").Value = 0.7323514225454915
$ws.Range("K16").Value = 0.3749417698318211
$ws.Range("L16").Value = 0.3177521848782874
$ws.Range("M16").Value = 0.1970830111043185
$ws.Range("O16").Value = 2.66634366643936
# Row 17
$ws.Range("B17").Value = 0.5745101600064686
$ws.Range("C17").Value = 0.3716660084953638
$ws.Range("D17").Value = 0.05109842788682073
$ws.Range("F17").Value = 0.7722842719497081
$ws.Range("G17").Value = 0.6099344706725347
$ws.Range("H17").Value = 0.7341452847657166
$ws.Range("K17").Value = 0.3600114860153951
$ws.Range("L17").Value = 0.3158301586307886
$ws.Range("M17").Value = 0.1932107801089771
$ws.Range("O17").Value = 2.671747765045012
# Row 18
$ws.Range("B18").Value = 0.5639576131237902
$ws.Range("C18").Value = 0.3722266913284749
$ws.Range("D18").Value = 0.05034556668457668
$ws.Range("F18").Value = 0.7722699877240302
$ws.Range("G18").Value = 0.6104454546748883
$ws.Range("H18").Value = 0.7352064064733455
$ws.Range("K18").Value = 0.3514148828560621
$ws.Range("L18").Value = 0.3147398682256721
$ws.Range("M18").Value = 0.1909896085441716
$ws.Range("O18").Value = 2.674997048928674
# Row 19
$ws.Range("B19").Value = 0.5603857007250781
$ws.Range("C19").Value = 0.3724180378631878
$ws.Range("D19").Value = 0.05009038453167847
$ws.Range("F19").Value = 0.7722726639316306
$ws.Range("G19").Value = 0.6106249774866086
$ws.Range("H19").Value = 0.7355707247862355
$ws.Range("K19").Value = 0.3485026747072482
$ws.Range("L19").Value = 0.314373329100988
$ws.Range("M19").Value = 0.1902385987115807
$ws.Range("O19").Value = 2.676121416639091
# Row 20
$ws.Range("B20").Value = 0.5764636686699589
$ws.Range("C20").Value = 0.3715629551844124
$ws.Range("D20").Value = 0.05123763406480464
$ws.Range("F20").Value = 0.7722904837658362
$ws.Range("G20").Value = 0.6098429926433155
$ws.Range("H20").Value = 0.7339512889181492
$ws.Range("K20").Value = 0.361601785876104
$ws.Range("L20").Value = 0.3160331881893086
$ws.Range("M20").Value = 0.193622362349636
$ws.Range("O20").Value = 2.671157897889316
# Row 21
$ws.Range("B21").Value = 0.6305498033690355
$ws.Range("C21").Value = 0.3687948067427662
$ws.Range("D21").Value = 0.05507284135791224
$ws.Range("F21").Value = 0.7728727609758366
$ws.Range("G21").Value = 0.6076688186892767
$ws.Range("H21").Value = 0.7288606301642062
$ws.Range("K21").Value = 0.40550306064182
$ws.Range("L21").Value = 0.3217916119865691
$ws.Range("M21").Value = 0.2050630984792718
$ws.Range("O21").Value = 2.656165596488108
# Row 22
$ws.Range("B22").Value = 0.6659440300457504
$ws.Range("C22").Value = 0.3670635948745691
$ws.Range("D22").Value = 0.05756471129180341
$ws.Range("F22").Value = 0.7736424209613801
$ws.Range("G22").Value = 0.6065858155967021
$ws.Range("H22").Value = 0.7257944918800803
$ws.Range("K22").Value = 0.4341103254242
$ws.Range("L22").Value = 0.325689953979051
$ws.Range("M22").Value = 0.2125930047197997
$ws.Range("O22").Value = 2.647622528671349
# Row 23
$ws.Range("B23").Value = 0.6470495077518024
$ws.Range("C23").Value = 0.3679804561370581
$ws.Range("D23").Value = 0.05623610795834111
$ws.Range("F23").Value = 0.7731960516148249
$ws.Range("G23").Value = 0.6071329014507114
$ws.Range("H23").Value = 0.7274070926789165
$ws.Range("K23").Value = 0.418850025388565
$ws.Range("L23").Value = 0.3235970206692826
$ws.Range("M23").Value = 0.2085693763276382
$ws.Range("O23").Value = 2.652067278882356
# Row 24
$ws.Range("B24").Value = 0.5755804845112777
$ws.Range("C24").Value = 0.3716095174596852
$ws.Range("D24").Value = 0.05117470504023203
$ws.Range("F24").Value = 0.7722875392578246
$ws.Range("G24").Value = 0.6098842310658696
$ws.Range("H24").Value = 0.7340389015917168
$ws.Range("K24").Value = 0.3608828517123754
$ws.Range("L24").Value = 0.3159413527065169
$ws.Range("M24").Value = 0.1934362703537715
$ws.Range("O24").Value = 2.671424133275565
# Row 25
$ws.Range("B25").Value = 0.4987096491629188
$ws.Range("C25").Value = 0.3758598570481304
$ws.Range("D25").Value = 0.04565353267349792
$ws.Range("F25").Value = 0.7729797046343663
$ws.Range("G25").Value = 0.6143017301046783
$ws.Range("H25").Value = 0.742314316386917
$ws.Range("K25").Value = 0.2980100089758082
$ws.Range("L25").Value = 0.3082651972506767
$ws.Range("M25").Value = 0.1773439588834833
$ws.Range("O25").Value = 2.697694030629805
